$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Reference No"
$ws.Range("B1").Value = "Customer"
$ws.Range("C1").Value = "Created At"
$ws.Range("D1").Value = "Service Cost"

# ---- Data rows: reference numbers, then customers (column-by-column), ----
# ---- matching the shared-string insertion order of the source report. ----
$ws.Range("A2").Value = "REF-20250731-001"
$ws.Range("A3").Value = "REF-20250731-002"
$ws.Range("A4").Value = "REF-20250731-003"

$ws.Range("B2").Value = "CRESTANKS LIMITED"
$ws.Range("B3").Value = "KYAGALANYI COFFEE LIMITED"
$ws.Range("B4").Value = "UGANDA REVENUE AUTHORITY"

# "Created At" column holds a literal text date string, not a real date
# serial. Stage it as a formula that evaluates to the literal text (so no
# date auto-parsing and no NumberFormat mutation ever happens), then copy
# just the VALUE (xlPasteValues) into the destination cells so the text
# lands as a plain shared string with no style attached to the cell and no
# stray number-format left behind in the style table.
$dateStage = $ws.Range("Z1")
$dateStage.Formula = '="2025-07-30"'
$dateStage.Copy()
$ws.Range("C2").PasteSpecial(-4163)
$ws.Range("C3").PasteSpecial(-4163)
$ws.Range("C4").PasteSpecial(-4163)
$dateStage.Clear()

$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 10
$ws.Range("D4").Value = 10

# ---- Header styling ----
# Build the combined look (bold font, thin box border, centered horizontally,
# top vertically) on a scratch cell and copy/paste the formats in a single
# shot, so only one merged cell style is minted rather than one per
# individually-set property.
$styleStage = $ws.Range("Z1")
$styleStage.Font.Bold = $true
$styleStage.Borders.LineStyle = 1
$styleStage.HorizontalAlignment = -4108
$styleStage.VerticalAlignment = -4160
$styleStage.Copy()

$headerRange = $ws.Range("A1:D1")
$headerRange.PasteSpecial(-4122)
$styleStage.Clear()
